$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A=M1, D=ECs
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Cdh1"
$ws.Range("C2").Value = "Igf1r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4197659999999999
$ws.Range("H2").Value = 1.259298
$ws.Range("I2").Value = 0.4104130959562529
$ws.Range("J2").Value = 0.4104130959562529
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.42467
$ws.Range("N2").Value = 37.27401
$ws.Range("O2").Value = 0.3327007343951245
$ws.Range("P2").Value = 0.3327007343951245
$ws.Range("Q2").Value = 5.215454027219999
$ws.Range("R2").Value = 46.93908624497999
$ws.Range("S2").Value = 0.136544738430022
$ws.Range("T2").Value = 0.136544738430022

# Row 3: A=M1, D=FAPs
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Cdh1"
$ws.Range("C3").Value = "Igf1r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4197659999999999
$ws.Range("H3").Value = 1.259298
$ws.Range("I3").Value = 0.4104130959562529
$ws.Range("J3").Value = 0.4104130959562529
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.45773566666667
$ws.Range("N3").Value = 37.373207
$ws.Range("O3").Value = 0.3335861479782027
$ws.Range("P3").Value = 0.3335861479782027
$ws.Range("Q3").Value = 5.229333869853999
$ws.Range("R3").Value = 47.064004828686
$ws.Range("S3").Value = 0.1369081237598549
$ws.Range("T3").Value = 0.1369081237598549

# Row 4: A=M1, D=M1
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Cdh1"
$ws.Range("C4").Value = "Igf1r"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4197659999999999
$ws.Range("H4").Value = 1.259298
$ws.Range("I4").Value = 0.4104130959562529
$ws.Range("J4").Value = 0.4104130959562529
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.546802333333333
$ws.Range("N4").Value = 10.640407
$ws.Range("O4").Value = 0.09497425211730705
$ws.Range("P4").Value = 0.09497425211730702
$ws.Range("Q4").Value = 1.488827028254
$ws.Range("R4").Value = 13.399443254286
$ws.Range("S4").Value = 0.03897867684759369
$ws.Range("T4").Value = 0.03897867684759368

# Row 5: A=M1, D=M2
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Cdh1"
$ws.Range("C5").Value = "Igf1r"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4197659999999999
$ws.Range("H5").Value = 1.259298
$ws.Range("I5").Value = 0.4104130959562529
$ws.Range("J5").Value = 0.4104130959562529
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.801381000000001
$ws.Range("N5").Value = 8.404143000000001
$ws.Range("O5").Value = 0.07501378435166073
$ws.Range("P5").Value = 0.07501378435166071
$ws.Range("Q5").Value = 1.175924496846
$ws.Range("R5").Value = 10.583320471614
$ws.Range("S5").Value = 0.0307866394751598
$ws.Range("T5").Value = 0.03078663947515979

# Row 6: A=M1, D=sCs
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Cdh1"
$ws.Range("C6").Value = "Igf1r"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4197659999999999
$ws.Range("H6").Value = 1.259298
$ws.Range("I6").Value = 0.4104130959562529
$ws.Range("J6").Value = 0.4104130959562529
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.114294
$ws.Range("N6").Value = 18.342882
$ws.Range("O6").Value = 0.1637250811577051
$ws.Range("P6").Value = 0.1637250811577051
$ws.Range("Q6").Value = 2.566572735204
$ws.Range("R6").Value = 23.099154616836
$ws.Range("S6").Value = 0.06719491744362252
$ws.Range("T6").Value = 0.0671949174436225

# Row 7: A=M2, D=ECs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cdh1"
$ws.Range("C7").Value = "Igf1r"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5921903333333333
$ws.Range("H7").Value = 1.776571
$ws.Range("I7").Value = 0.5789956025468922
$ws.Range("J7").Value = 0.5789956025468922
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 12.42467
$ws.Range("N7").Value = 37.27401
$ws.Range("O7").Value = 0.3327007343951245
$ws.Range("P7").Value = 0.3327007343951245
$ws.Range("Q7").Value = 7.357769468856665
$ws.Range("R7").Value = 66.21992521970999
$ws.Range("S7").Value = 0.1926322621788986
$ws.Range("T7").Value = 0.1926322621788986

# Row 8: A=M2, D=FAPs
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cdh1"
$ws.Range("C8").Value = "Igf1r"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5921903333333333
$ws.Range("H8").Value = 1.776571
$ws.Range("I8").Value = 0.5789956025468922
$ws.Range("J8").Value = 0.5789956025468922
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.45773566666667
$ws.Range("N8").Value = 37.373207
$ws.Range("O8").Value = 0.3335861479782027
$ws.Range("P8").Value = 0.3335861479782027
$ws.Range("Q8").Value = 7.377350637021888
$ws.Range("R8").Value = 66.39615573319699
$ws.Range("S8").Value = 0.1931449127499362
$ws.Range("T8").Value = 0.1931449127499362

# Row 9: A=M2, D=M1
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cdh1"
$ws.Range("C9").Value = "Igf1r"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5921903333333333
$ws.Range("H9").Value = 1.776571
$ws.Range("I9").Value = 0.5789956025468922
$ws.Range("J9").Value = 0.5789956025468922
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.546802333333333
$ws.Range("N9").Value = 10.640407
$ws.Range("O9").Value = 0.09497425211730705
$ws.Range("P9").Value = 0.09497425211730702
$ws.Range("Q9").Value = 2.100382056044111
$ws.Range("R9").Value = 18.903438504397
$ws.Range("S9").Value = 0.05498967433110065
$ws.Range("T9").Value = 0.05498967433110063

# Row 10: A=M2, D=M2
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Cdh1"
$ws.Range("C10").Value = "Igf1r"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5921903333333333
$ws.Range("H10").Value = 1.776571
$ws.Range("I10").Value = 0.5789956025468922
$ws.Range("J10").Value = 0.5789956025468922
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.801381000000001
$ws.Range("N10").Value = 8.404143000000001
$ws.Range("O10").Value = 0.07501378435166073
$ws.Range("P10").Value = 0.07501378435166071
$ws.Range("Q10").Value = 1.658950748183667
$ws.Range("R10").Value = 14.930556733653
$ws.Range("S10").Value = 0.04343265127001244
$ws.Range("T10").Value = 0.04343265127001243

# Row 11: A=M2, D=sCs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Cdh1"
$ws.Range("C11").Value = "Igf1r"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5921903333333333
$ws.Range("H11").Value = 1.776571
$ws.Range("I11").Value = 0.5789956025468922
$ws.Range("J11").Value = 0.5789956025468922
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.114294
$ws.Range("N11").Value = 18.342882
$ws.Range("O11").Value = 0.1637250811577051
$ws.Range("P11").Value = 0.1637250811577051
$ws.Range("Q11").Value = 3.620825801958
$ws.Range("R11").Value = 32.58743221762199
$ws.Range("S11").Value = 0.09479610201694429
$ws.Range("T11").Value = 0.09479610201694427

# Row 12: A=sCs, D=ECs
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Cdh1"
$ws.Range("C12").Value = "Igf1r"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.01083266666666667
$ws.Range("H12").Value = 0.032498
$ws.Range("I12").Value = 0.01059130149685484
$ws.Range("J12").Value = 0.01059130149685484
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 12.42467
$ws.Range("N12").Value = 37.27401
$ws.Range("O12").Value = 0.3327007343951245
$ws.Range("P12").Value = 0.3327007343951245
$ws.Range("Q12").Value = 0.1345923085533333
$ws.Range("R12").Value = 1.21133077698
$ws.Range("S12").Value = 0.003523733786203787
$ws.Range("T12").Value = 0.003523733786203787

# Row 13: A=sCs, D=FAPs
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Cdh1"
$ws.Range("C13").Value = "Igf1r"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.01083266666666667
$ws.Range("H13").Value = 0.032498
$ws.Range("I13").Value = 0.01059130149685484
$ws.Range("J13").Value = 0.01059130149685484
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 12.45773566666667
$ws.Range("N13").Value = 37.373207
$ws.Range("O13").Value = 0.3335861479782027
$ws.Range("P13").Value = 0.3335861479782027
$ws.Range("Q13").Value = 0.1349504978984444
$ws.Range("R13").Value = 1.214554481086
$ws.Range("S13").Value = 0.003533111468411579
$ws.Range("T13").Value = 0.003533111468411578

# Row 14: A=sCs, D=M1
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Cdh1"
$ws.Range("C14").Value = "Igf1r"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.01083266666666667
$ws.Range("H14").Value = 0.032498
$ws.Range("I14").Value = 0.01059130149685484
$ws.Range("J14").Value = 0.01059130149685484
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.546802333333333
$ws.Range("N14").Value = 10.640407
$ws.Range("O14").Value = 0.09497425211730705
$ws.Range("P14").Value = 0.09497425211730702
$ws.Range("Q14").Value = 0.03842132740955555
$ws.Range("R14").Value = 0.345791946686
$ws.Range("S14").Value = 0.001005900938612703
$ws.Range("T14").Value = 0.001005900938612703

# Row 15: A=sCs, D=M2
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Cdh1"
$ws.Range("C15").Value = "Igf1r"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.01083266666666667
$ws.Range("H15").Value = 0.032498
$ws.Range("I15").Value = 0.01059130149685484
$ws.Range("J15").Value = 0.01059130149685484
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.801381000000001
$ws.Range("N15").Value = 8.404143000000001
$ws.Range("O15").Value = 0.07501378435166073
$ws.Range("P15").Value = 0.07501378435166071
$ws.Range("Q15").Value = 0.03034642657933334
$ws.Range("R15").Value = 0.273117839214
$ws.Range("S15").Value = 0.0007944936064884906
$ws.Range("T15").Value = 0.0007944936064884904

# Row 16: A=sCs, D=sCs
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Cdh1"
$ws.Range("C16").Value = "Igf1r"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.01083266666666667
$ws.Range("H16").Value = 0.032498
$ws.Range("I16").Value = 0.01059130149685484
$ws.Range("J16").Value = 0.01059130149685484
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 6.114294
$ws.Range("N16").Value = 18.342882
$ws.Range("O16").Value = 0.1637250811577051
$ws.Range("P16").Value = 0.1637250811577051
$ws.Range("Q16").Value = 0.06623410880399999
$ws.Range("R16").Value = 0.596106979236
$ws.Range("S16").Value = 0.001734061697138282
$ws.Range("T16").Value = 0.001734061697138282
